$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    @("42÷9=4, 6", "35÷9=3, 8"),
    @("68÷3=22, 2", "19÷7=2, 5"),
    @("92÷7=13, 1", "94÷5=18, 4"),
    @("16÷4=4, 0", "45÷6=7, 3"),
    @("20÷7=2, 6", "29÷6=4, 5"),
    @("77÷4=19, 1", "63÷5=12, 3"),
    @("47÷5=9, 2", "67÷6=11, 1"),
    @("44÷6=7, 2", "48÷2=24, 0"),
    @("84÷8=10, 4", "25÷8=3, 1"),
    @("99÷8=12, 3", "54÷3=18, 0"),
    @("55÷9=6, 1", "88÷6=14, 4"),
    @("29÷6=4, 5", "95÷2=47, 1"),
    @("93÷7=13, 2", "33÷6=5, 3"),
    @("68÷3=22, 2", "18÷8=2, 2"),
    @("92÷3=30, 2", "88÷3=29, 1"),
    @("56÷2=28, 0", "31÷3=10, 1"),
    @("95÷9=10, 5", "32÷2=16, 0"),
    @("74÷5=14, 4", "90÷5=18, 0"),
    @("16÷7=2, 2", "96÷4=24, 0"),
    @("77÷2=38, 1", "89÷2=44, 1"),
    @("49÷5=9, 4", "50÷9=5, 5"),
    @("66÷2=33, 0", "24÷7=3, 3"),
    @("76÷4=19, 0", "27÷9=3, 0"),
    @("24÷5=4, 4", "63÷3=21, 0"),
    @("18÷8=2, 2", "67÷2=33, 1")
)

$rows = @(1, 5, 9, 13, 17)
$idx = 0

foreach ($r in $rows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $pair = $values[$idx]
        $old = $pair[0]
        $new = $pair[1]
        $cellRange = $cell.Range
        # Trim the end-of-cell marker off the range so we only overwrite the
        # visible text, then replace it with the new equation.
        $cellRange.MoveEnd(1, -1) | Out-Null
        $current = $cellRange.Text
        if ($current -ne $old) {
            Write-Output ("Warning: row $r col $c expected '" + $old + "' but found '" + $current + "'")
        }
        $cellRange.Text = $new
        $idx = $idx + 1
    }
}
